$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.268.53"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").Value = "3.174.98"
$ws.Range("E3").Value = "  -8.24%  "
$ws.Range("D5").Value = "'565.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.07%  "
$ws.Range("D6").Value = "'169.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "3.173.88"
$ws.Range("E9").Value = "  -8.23%  "
$ws.Range("E10").Value = "  -6.50%  "
$ws.Range("D11").Value = "'6.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.67%  "
$ws.Range("E12").Value = "  -4.82%  "
$ws.Range("D13").Value = "3.725.06"
$ws.Range("E13").Value = "  -8.22%  "
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "'27.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.33%  "
$ws.Range("D16").Value = "64.276.97"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("E17").Value = "  -5.28%  "
$ws.Range("D18").Value = "3.171.93"
$ws.Range("E18").Value = "  -8.20%  "
$ws.Range("E19").Value = "  -3.93%  "
$ws.Range("D20").Value = "'13.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.04%  "
$ws.Range("D21").Value = "'353.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.53%  "
$ws.Range("D22").Value = "'7.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.84%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'69.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.64%  "
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("E26").Value = "  -5.85%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").Value = "'0.176"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.38%  "
$ws.Range("E32").Value = "  -5.31%  "
$ws.Range("D33").Value = "'22.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.06%  "
$ws.Range("E34").Value = "  -5.99%  "
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("E36").Value = "  -7.05%  "
$ws.Range("D37").Value = "'153.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("E38").Value = "  -7.65%  "
$ws.Range("D39").Value = "'25.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.66%  "
$ws.Range("D40").Value = "'2.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -6.19%  "
$ws.Range("D42").Value = "2.625.88"
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("E43").Value = "  -7.42%  "
$ws.Range("D44").Value = "'6.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.80%  "
$ws.Range("D45").Value = "'39.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("D47").Value = "'23.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("D48").Value = "'321.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.38%  "
$ws.Range("E49").Value = "  -7.54%  "
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.02%  "
